$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.143.34'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = '2.469.39'
$ws.Range('E3').Value = '  -1.09%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.21'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.04'
$ws.Range('E6').Value = '  -3.80%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -2.30%  '
$ws.Range('D9').Value = '2.469.44'
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('E10').Value = '  -4.50%  '
$ws.Range('E11').Value = '  -1.00%  '
$ws.Range('E12').Value = '  -3.34%  '
$ws.Range('E13').Value = '  -2.83%  '
$ws.Range('E14').Value = '  -0.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.47'
$ws.Range('E15').Value = '  -3.12%  '
$ws.Range('D16').Value = '67.052.93'
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('E17').Value = '  -4.52%  '
$ws.Range('D18').Value = '2.463.88'
$ws.Range('E18').Value = '  -1.45%  '
$ws.Range('E19').Value = '  -2.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.65'
$ws.Range('E20').Value = '  -4.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '356.35'
$ws.Range('E21').Value = '  -2.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.03'
$ws.Range('E22').Value = '  -2.63%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.57'
$ws.Range('E24').Value = '  -2.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.24'
$ws.Range('E25').Value = '  -7.20%  '
$ws.Range('E26').Value = '  -7.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.99'
$ws.Range('E27').Value = '  -9.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('D29').Value = '2.595.25'
$ws.Range('E29').Value = '  -0.51%  '
$ws.Range('E30').Value = '  -6.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '512.54'
$ws.Range('E31').Value = '  -4.08%  '
$ws.Range('E32').Value = '  -5.41%  '
$ws.Range('E33').Value = '  -4.61%  '
$ws.Range('E34').Value = '  -5.94%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.119'
$ws.Range('E36').Value = '  -7.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.74'
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.58'
$ws.Range('E38').Value = '  -0.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.46'
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('E40').Value = '  -6.09%  '
$ws.Range('E41').Value = '  -6.70%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.326'
$ws.Range('E42').Value = '  -6.42%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.79'
$ws.Range('E43').Value = '  -6.56%  '
$ws.Range('E44').Value = '  -2.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.31'
$ws.Range('E45').Value = '  -7.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '141.50'
$ws.Range('E46').Value = '  -2.25%  '
$ws.Range('E47').Value = '  -5.55%  '
$ws.Range('E48').Value = '  -5.83%  '
$ws.Range('E49').Value = '  -7.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.59'
$ws.Range('E50').Value = '  -6.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0732'
$ws.Range('E51').Value = '  -2.33%  '
